$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-06 02:28:24"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
